$d = $word.ActiveDocument

# The document contains six paragraphs, each holding the literal text
# "<id>p100r_N</id>" split across three separate runs (the "<id>" open
# tag and "</id>" close tag in Courier New / color 7f6000, and the
# "p100r_N" body in plain black). The edit merges each trio of runs
# into a single run (keeping the first run's formatting) whose text is
# the full "<id>p100r_N</id>" string.
#
# Doing a Find/Replace over the full span with MatchWildcards off makes
# Word collapse the matched runs into one run carrying the formatting
# of the first run in the match, which reproduces the target XML.

$ids = 1..6

foreach ($n in $ids) {
    $old = "<id>p100r_$n</id>"
    $new = "<id>p100r_$n</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
